$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Update the username/email value used in A2:A4 (shared string content change)
$ws.Range("A2").Value = "likitha.lokesh@slalom.com"
$ws.Range("A3").Value = "likitha.lokesh@slalom.com"
$ws.Range("A4").Value = "likitha.lokesh@slalom.com"

# Move the active selection from B4 to B5
$ws.Range("B5").Select()
